$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values regenerated for column G (rows 2-33),
# replacing the old "Strike#" derived figures.
$newK = @{
    2  = 2
    3  = 3
    4  = 5
    5  = 3
    6  = 1
    7  = 10
    8  = 4
    9  = 9
    10 = 8
    11 = 7
    12 = 2
    13 = 3
    14 = 4
    15 = 4
    16 = 5
    17 = 6
    18 = 2
    19 = 3
    20 = 2
    21 = 3
    22 = 7
    23 = 7
    24 = 1
    25 = 5
    26 = 3
    27 = 9
    28 = 4
    29 = 2
    30 = 3
    31 = 6
    32 = 2
    33 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
